$p = $ppt.ActivePresentation

# Find the slide whose title contains the typo "problema" and fix it to "problem",
# matching the run-split produced by editing just the word (plus its trailing
# space) in place, leaving the surrounding runs ("the ... " / " (cont.)") intact.
foreach ($s in $p.Slides) {
    foreach ($sh in $s.Shapes) {
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            $tr = $sh.TextFrame.TextRange
            $full = $tr.Text
            $idx = $full.IndexOf("problema ")
            if ($idx -ge 0) {
                $sub = $tr.Characters($idx + 1, 9)
                $sub.Text = "problem "
            }
        }
    }
}
